$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings (e.g. "1.004", "7.960") keep their exact original formatting
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply all the updated cell values from the crypto price refresh.
$ws.Range("D2").Value = '25.080.53'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").Value = '1.705.89'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '316.53'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.4004'
$ws.Range("E7").Value = '  +2.43%  '
$ws.Range("D8").Value = '0.4037'
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("D9").Value = '1.470'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '52.87'
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '0.08834'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '25.96'
$ws.Range("E13").Value = '  -2.72%  '
$ws.Range("D14").Value = '7.493'
$ws.Range("D15").Value = '0.00001355'
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '7.960'
$ws.Range("E16").Value = '  -3.58%  '
$ws.Range("D17").Value = '1.706.65'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '96.30'
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("D19").Value = '0.07196'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '20.78'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '7.291'
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = '14.35'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").Value = '25.108.89'
$ws.Range("E24").Value = '  +2.40%  '
$ws.Range("D25").Value = '2.392'
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("D26").Value = '2.937'
$ws.Range("E26").Value = '  -2.44%  '
$ws.Range("D27").Value = '23.63'
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("D28").Value = '6.213'
$ws.Range("E28").Value = '  +14.63%  '
$ws.Range("D29").Value = '162.99'
$ws.Range("E29").Value = '  -3.45%  '
$ws.Range("D30").Value = '151.44'
$ws.Range("E30").Value = '  +3.97%  '
$ws.Range("D31").Value = '8.341'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").Value = '2.669'
$ws.Range("E32").Value = '  +22.22%  '
$ws.Range("D33").Value = '1.892.46'
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D34").Value = '0.08607'
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("D35").Value = '0.03161'
$ws.Range("E35").Value = '  +3.43%  '
$ws.Range("D36").Value = '1.045'
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = '7.196'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '0.2911'
$ws.Range("E38").Value = '  +3.57%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = '0.09754'
$ws.Range("E39").Value = '  +6.20%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '10.99'
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").Value = '0.8299'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("E42").Value = '  -1.22%  '
$ws.Range("D43").Value = '1.477'
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Value = '17.13'
$ws.Range("E44").Value = '  -2.56%  '
$ws.Range("D45").Value = '2.683'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").Value = '0.7421'
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("B47").Value = 'Flow'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D47").Value = '1.437'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.09102'
$ws.Range("E48").Value = '  +10.95%  '
$ws.Range("D49").Value = '4.251'
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").Value = '1.002'
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").Value = '140.29'
$ws.Range("E51").Value = '  -0.03%  '

# Restore the default "Normal" style on column D so the cells end up
# exactly as they were stylistically (no explicit style index), matching
# the rest of the untouched data cells.
$ws.Range("D2:D51").Style = "Normal"
